$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "20.551.35"
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +1.30%  "

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "1.471.25"
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +2.05%  "

$ws.Cells.Item(4, 5).Value = "  +0.45%  "

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "0.9576"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +4.80%  "

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "276.95"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +0.73%  "

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.3560"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -1.93%  "

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.3069"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -0.30%  "

$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "1.093"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +6.70%  "

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "39.42"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +0.39%  "

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.06633"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +1.77%  "

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +0.56%  "

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "5.459"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +1.97%  "

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "18.08"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +3.01%  "

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "6.174"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +2.13%  "

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "0.9580"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +3.02%  "

$ws.Cells.Item(17, 5).Value = "  +0.95%  "

$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "1.469.87"
$c.Style = "Normal"

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "0.05961"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +5.98%  "

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "68.95"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +1.99%  "

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "5.482"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +1.13%  "

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "14.52"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +1.96%  "

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "11.27"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +4.08%  "

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "2.267"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +1.47%  "

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "20.570.62"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +1.42%  "

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "145.06"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +5.34%  "

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "2.086"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -0.20%  "

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "17.12"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +1.17%  "

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "1.628.83"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +2.38%  "

$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "114.00"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +3.32%  "

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "3.835"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -2.34%  "

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "4.930"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +1.40%  "

$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "0.07927"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +3.36%  "

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "0.7953"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -1.19%  "

$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "1.241"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +9.66%  "

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "1.436"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -1.63%  "

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "0.05733"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -1.42%  "

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "4.713"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +0.67%  "

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "0.02027"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +1.90%  "

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "0.9580"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +3.41%  "

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "10.32"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +1.20%  "

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "0.1858"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +0.23%  "

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "7.272"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +2.20%  "

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "0.5245"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +0.54%  "

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "3.511"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +0.74%  "

$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "12.10"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +1.95%  "

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "118.56"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +1.30%  "

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "0.5179"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +1.42%  "

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "1.800"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +3.87%  "

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "0.06434"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +0.31%  "

$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "0.9895"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +1.35%  "
